$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.785.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.814.39'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.14%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '276.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -8.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5085'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3515'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.24'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06671'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.01'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8312'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07917'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.812.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.084'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.81%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9991'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.36%  '
$ws.Range("E19").Value = '  -5.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.0000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '25.823.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.724'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.081'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.83%  '
$ws.Range("E25").Value = '  -3.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.172'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.48%  '
$ws.Range("E27").Value = '  -3.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '109.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.323'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.239'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08835'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04868'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7352'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.136'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.882'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.149'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9994'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5223'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -11.81%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.311'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -11.63%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01843'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9586'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '112.84'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.195'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.077'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9997'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4592'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1363'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.78%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.43'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.26%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.250'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.503'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.29%  '
